$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1 (20:22 -> 20:52)
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 20:52"

# Row 23 = Irlanda: update new-case figures
$ws.Range("B23").Value = 18184
$ws.Range("C23").Value = 577
$ws.Range("E23").Value = 7937
$ws.Range("F23").Value = 142
$ws.Range("G23").Value = 220
$ws.Range("H23").Value = 1014

# Sudafrica overtakes Egipto in the ranking: Sudafrica now sits in row 53
# (above Egipto, which moves to row 54). Egipto keeps its previous figures,
# Sudafrica gets updated figures.
$ws.Range("A53").Value = "Sudafrica"
$ws.Range("B53").Value = 4220
$ws.Range("C53").Value = 267
$ws.Range("D53").Value = 1473
$ws.Range("E53").Value = 2668
$ws.Range("F53").Value = 36
$ws.Range("G53").Value = 4
$ws.Range("H53").Value = 79

$ws.Range("A54").Value = "Egipto"
$ws.Range("B54").Value = 4092
$ws.Range("C54").Value = 201
$ws.Range("D54").Value = 1075
$ws.Range("E54").Value = 2723
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 7
$ws.Range("H54").Value = 294

# Row 120 = Venezuela: update new-case figures
$ws.Range("B120").Value = 318
$ws.Range("C120").Value = 20
$ws.Range("D120").Value = 128
$ws.Range("E120").Value = 180
